$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D (Price) and E (Volume) that need updated values.
# They hold text-like values (e.g. "325.45", "2.58%"), so force the
# affected ranges to Text format before writing, otherwise Excel would
# coerce them into numeric values (changing cell type/formatting).
# (Applied as two contiguous blocks since a multi-area Range only
# formats its first area reliably.)
$ws.Range("D2:E27").NumberFormat = "@"
$ws.Range("D39:E51").NumberFormat = "@"

$ws.Range("D2").Value = "325.45"
$ws.Range("E2").Value = "2.58%"
$ws.Range("D3").Value = "39.99"
$ws.Range("E3").Value = "5.46%"
$ws.Range("D4").Value = "5.863"
$ws.Range("E4").Value = "12.99%"
$ws.Range("D5").Value = "0.08000"
$ws.Range("E5").Value = "-0.49%"
$ws.Range("D6").Value = "4.581"
$ws.Range("E6").Value = "1.89%"
$ws.Range("D7").Value = "8.711"
$ws.Range("E7").Value = "2.35%"
$ws.Range("D8").Value = "1.914"
$ws.Range("E8").Value = "-1.73%"
$ws.Range("E9").Value = "-1.02%"
$ws.Range("D10").Value = "0.9397"
$ws.Range("E10").Value = "-0.10%"
$ws.Range("D11").Value = "0.1249"
$ws.Range("E11").Value = "-3.15%"
$ws.Range("D12").Value = "0.1964"
$ws.Range("E12").Value = "1.47%"
$ws.Range("D13").Value = "8.837"
$ws.Range("E13").Value = "35.04%"
$ws.Range("D14").Value = "0.09184"
$ws.Range("E14").Value = "2.00%"
$ws.Range("D15").Value = "0.03582"
$ws.Range("E15").Value = "5.59%"
$ws.Range("D16").Value = "0.09623"
$ws.Range("E16").Value = "0.80%"
$ws.Range("D17").Value = "0.001313"
$ws.Range("E17").Value = "-5.58%"
$ws.Range("E18").Value = "1.42%"
$ws.Range("E19").Value = "-0.71%"
$ws.Range("D20").Value = "0.3525"
$ws.Range("E20").Value = "0.06%"
$ws.Range("D21").Value = "0.1432"
$ws.Range("E21").Value = "8.84%"
$ws.Range("E22").Value = "-0.37%"
$ws.Range("D23").Value = "0.04443"
$ws.Range("E23").Value = "1.34%"
$ws.Range("D24").Value = "0.001260"
$ws.Range("E24").Value = "2.44%"
$ws.Range("D25").Value = "0.004318"
$ws.Range("E25").Value = "1.11%"
$ws.Range("D26").Value = "0.0001145"
$ws.Range("E26").Value = "-13.91%"
$ws.Range("E27").Value = "0.09%"
$ws.Range("D39").Value = "0.02426"
$ws.Range("E39").Value = "2.92%"
$ws.Range("D40").Value = "0.05260"
$ws.Range("E40").Value = "1.97%"
$ws.Range("D41").Value = "0.007507"
$ws.Range("E41").Value = "-2.36%"
$ws.Range("E42").Value = "0.92%"
$ws.Range("D43").Value = "0.008683"
$ws.Range("E43").Value = "-0.03%"
$ws.Range("D44").Value = "0.002106"
$ws.Range("E44").Value = "-0.19%"
$ws.Range("D45").Value = "0.01059"
$ws.Range("E45").Value = "19.63%"
$ws.Range("D46").Value = "0.00006843"
$ws.Range("E46").Value = "5.52%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.27%"
$ws.Range("D48").Value = "0.002886"
$ws.Range("E48").Value = "0.90%"
$ws.Range("D49").Value = "0.001424"
$ws.Range("E49").Value = "-15.71%"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").Value = "0.27%"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.27%"

# Restore the default (Normal) cell style so formatting matches the
# original workbook - only the cell contents should change.
$ws.Range("D2:E27").Style = "Normal"
$ws.Range("D39:E51").Style = "Normal"
